$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1624.75
$ws.Range("I6").Value = 1624.75
$ws.Range("K6").Value = 4874.25
$ws.Range("M6").Value = -4762.25
$ws.Range("H8").Value = 62530.375
$ws.Range("I8").Value = 62530.375
$ws.Range("K8").Value = 187591.125
$ws.Range("M8").Value = -187452.125
$ws.Range("H17").Value = 5023089.5
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 5262260.5
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 15786781.5
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -15787117.5
$ws.Range("H31").Value = 601.5
$ws.Range("I31").Value = 601.5
$ws.Range("K31").Value = 1804.5
$ws.Range("M31").Value = -1574.5
$ws.Range("H112").Value = 1845.3125
$ws.Range("I112").Value = 1066.6666
$ws.Range("J112").Value = 1925.862
$ws.Range("K112").Value = 3199.9998
$ws.Range("L112").Value = 5777.586
$ws.Range("M112").Value = -2091.9998
$ws.Range("N112").Value = -7993.586
$ws.Range("H137").Value = 7114.8125
$ws.Range("I137").Value = 774.5714
$ws.Range("K137").Value = 2323.7142
$ws.Range("M137").Value = 226.2857999999997
$ws.Range("H138").Value = 6065716
$ws.Range("I138").Value = 19611490
$ws.Range("J138").Value = 5764.4473
$ws.Range("K138").Value = 58834470
$ws.Range("L138").Value = 17293.3419
$ws.Range("M138").Value = -58829330
$ws.Range("N138").Value = -27573.3419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2066.889
$ws.Range("I61").Value = 1936.8292
$ws.Range("J61").Value = 3400
$ws.Range("K61").Value = 1936.8292
$ws.Range("L61").Value = 3400
$ws.Range("M61").Value = -1724.8292
$ws.Range("N61").Value = -3824
$ws.Range("H74").Value = 4040.4849
$ws.Range("I74").Value = 952.069
$ws.Range("J74").Value = 26431.5
$ws.Range("K74").Value = 952.069
$ws.Range("L74").Value = 26431.5
$ws.Range("M74").Value = -78.06899999999996
$ws.Range("N74").Value = -28179.5
$ws.Range("H77").Value = 4040.4849
$ws.Range("I77").Value = 952.069
$ws.Range("J77").Value = 26431.5
$ws.Range("K77").Value = 4760.344999999999
$ws.Range("L77").Value = 132157.5
$ws.Range("M77").Value = -392.3449999999993
$ws.Range("N77").Value = -140893.5
$ws.Range("H132").Value = 1603.2941
$ws.Range("I132").Value = 1255.5435
$ws.Range("J132").Value = 4802.6
$ws.Range("K132").Value = 3766.6305
$ws.Range("L132").Value = 14407.8
$ws.Range("M132").Value = -1236.6305
$ws.Range("N132").Value = -19467.8
$ws.Range("H136").Value = 2066.889
$ws.Range("I136").Value = 1936.8292
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 5810.487599999999
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = -3260.487599999999
$ws.Range("N136").Value = -15300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2413.4043
$ws.Range("I134").Value = 1788.1904
$ws.Range("J134").Value = 7665.2
$ws.Range("K134").Value = 5364.5712
$ws.Range("L134").Value = 22995.6
$ws.Range("M134").Value = -2829.5712
$ws.Range("N134").Value = -28065.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4047.663
$ws.Range("I31").Value = 1074.125
$ws.Range("J31").Value = 7291.523
$ws.Range("K31").Value = 1074.125
$ws.Range("L31").Value = 7291.523
$ws.Range("M31").Value = -779.125
$ws.Range("N31").Value = -7881.523
$ws.Range("H34").Value = 4047.663
$ws.Range("I34").Value = 1074.125
$ws.Range("J34").Value = 7291.523
$ws.Range("K34").Value = 1074.125
$ws.Range("L34").Value = 7291.523
$ws.Range("M34").Value = -872.125
$ws.Range("N34").Value = -7695.523
$ws.Range("H132").Value = 3858.5417
$ws.Range("I132").Value = 1721.4865
$ws.Range("J132").Value = 11046.818
$ws.Range("K132").Value = 5164.4595
$ws.Range("L132").Value = 33140.454
$ws.Range("M132").Value = -2634.4595
$ws.Range("N132").Value = -38200.454
$ws.Range("H141").Value = 36512
$ws.Range("I141").Value = 14500
$ws.Range("J141").Value = 39263.5
$ws.Range("K141").Value = 14500
$ws.Range("L141").Value = 39263.5
$ws.Range("M141").Value = -9320
$ws.Range("N141").Value = -49623.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 257.9
$ws.Range("I6").Value = 197.375
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 592.125
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -479.125
$ws.Range("N6").Value = -1726
$ws.Range("H11").Value = 238
$ws.Range("I11").Value = 208.88889
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 626.6666700000001
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = -486.6666700000001
$ws.Range("N11").Value = -1780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 21833.334
$ws.Range("J74").Value = 21833.334
$ws.Range("L74").Value = 21833.334
$ws.Range("N74").Value = -23705.334
$ws.Range("H77").Value = 21833.334
$ws.Range("J77").Value = 21833.334
$ws.Range("L77").Value = 65500.00199999999
$ws.Range("N77").Value = -74860.00199999999
$ws.Range("H118").Value = 19997.5
$ws.Range("J118").Value = 19997.5
$ws.Range("L118").Value = 19997.5
$ws.Range("N118").Value = -23311.5
$ws.Range("H121").Value = 45000
$ws.Range("J121").Value = 45000
$ws.Range("L121").Value = 45000
$ws.Range("N121").Value = -48494
$ws.Range("H126").Value = 4007.1538
$ws.Range("I126").Value = 2879.6
$ws.Range("K126").Value = 8638.799999999999
$ws.Range("M126").Value = -6168.799999999999
$ws.Range("H132").Value = 4783.757
$ws.Range("I132").Value = 5158.0967
$ws.Range("J132").Value = 2849.6667
$ws.Range("K132").Value = 15474.2901
$ws.Range("L132").Value = 8549.000100000001
$ws.Range("M132").Value = -12944.2901
$ws.Range("N132").Value = -13609.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8035.5
$ws.Range("I7").Value = 6331.75
$ws.Range("J7").Value = 10080
$ws.Range("K7").Value = 6331.75
$ws.Range("L7").Value = 10080
$ws.Range("M7").Value = -6219.75
$ws.Range("N7").Value = -10304
$ws.Range("H81").Value = 49828.6
$ws.Range("J81").Value = 49828.6
$ws.Range("L81").Value = 49828.6
$ws.Range("N81").Value = -51824.6
$ws.Range("H84").Value = 49828.6
$ws.Range("J84").Value = 49828.6
$ws.Range("L84").Value = 149485.8
$ws.Range("N84").Value = -159469.8
$ws.Range("H126").Value = 8035.5
$ws.Range("I126").Value = 6331.75
$ws.Range("J126").Value = 10080
$ws.Range("K126").Value = 18995.25
$ws.Range("L126").Value = 30240
$ws.Range("M126").Value = -16525.25
$ws.Range("N126").Value = -35180
$ws.Range("H136").Value = 6322.4375
$ws.Range("I136").Value = 2843.5454
$ws.Range("J136").Value = 13976
$ws.Range("K136").Value = 8530.636200000001
$ws.Range("L136").Value = 41928
$ws.Range("M136").Value = -5980.636200000001
$ws.Range("N136").Value = -47028
